# Update pin list for SD, EEPROM, Flash, RJ45 (and FPGA/Part/Alt/pin-list columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new "FPGA" divider row above each "DDR3 SODIMM" block ---
# First block (under "FBG484 Package" header on row 1): new row becomes row 2.
$ws.Rows.Item(2).Insert()
# Second block (under "FBG676 Package" header, now shifted to row 12): new row becomes row 13.
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(2, 1).Value = "FPGA"
$ws.Cells.Item(13, 1).Value = "FPGA"
# The inserted rows pick up the bold header formatting of the row above;
# these sub-headers are regular (non-bold) text, so clear that back out.
$ws.Cells.Item(2, 1).Font.Bold = $false
$ws.Cells.Item(13, 1).Font.Bold = $false

# --- Alternate part for the Quad SPI Flash rows ---
$ws.Cells.Item(4, 3).Value = "MT25Q"
$ws.Cells.Item(15, 3).Value = "MT25Q"

# --- New header columns on row 1 ---
$ws.Cells.Item(1, 2).Value = "Part"
$ws.Cells.Item(1, 3).Value = "Alt"
$ws.Cells.Item(1, 4).Value = "On pin list?"
$ws.Cells.Item(1, 5).Value = "On schematic?"

# --- "On pin list?" marks for the first (FBG484) block ---
$ws.Cells.Item(3, 4).Value = "X"
$ws.Cells.Item(4, 4).Value = "X"
$ws.Cells.Item(5, 4).Value = "X"
$ws.Cells.Item(6, 4).Value = "X"
$ws.Cells.Item(7, 4).Value = "X"
$ws.Cells.Item(8, 4).Value = "X"
$ws.Cells.Item(10, 4).Value = "X"

# --- Selection left where the author finished editing ---
$ws.Range("D3").Select()
